$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a weekly price log: a new week's worth of rows (3 rows,
# one per "Calidad" grade: Especial / Primera / Segunda) is published on
# top, and every existing row shifts down by three positions.
# Insert three blank rows at the top of the data block (row 2) so the
# existing rows 2..17 become rows 5..20.
$ws.Rows.Item(2).EntireRow.Insert()
$ws.Rows.Item(2).EntireRow.Insert()
$ws.Rows.Item(2).EntireRow.Insert()

# The freshly inserted rows inherit the bold/bordered header style from
# row 1. Reset that so the new data rows look like ordinary data rows.
$ws.Range("A2:T4").ClearFormats()

# New rows for the week of 2022-03-10 (serial 44630).
$newRows = @(
    @{ Row = 2; Calidad = "Especial"; Volumen = 300; Min = 15000; Max = 16000; Prom = 15500; KgKilo = 861 },
    @{ Row = 3; Calidad = "Primera";  Volumen = 300; Min = 12000; Max = 13000; Prom = 12500; KgKilo = 694 },
    @{ Row = 4; Calidad = "Segunda";  Volumen = 240; Min = 9000;  Max = 10000; Prom = 9500;  KgKilo = 528 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 4).Value = 44630
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107011
    $ws.Cells.Item($row, 10).Value = "Tuna"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = "$/caja 18 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.KgKilo
    $ws.Cells.Item($row, 20).Value = 18
}
